# Update res_bus/vm_pu.xlsx values for Case_2_17 ("case with 380 kV done").
# The voltage setpoint at bus column B moves from 1.05 p.u. to 1.02 p.u.
# for every time-step row (2-25), and the resulting load-flow voltage
# magnitudes in columns C:F and I:N are updated to the newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035893676681444
$ws.Range("D2").Value = 1.040614110272484
$ws.Range("E2").Value = 1.055220103135679
$ws.Range("F2").Value = 1.061537761288126
$ws.Range("I2").Value = 1.041648635642988
$ws.Range("J2").Value = 1.041004963729284
$ws.Range("K2").Value = 1.043395959668502
$ws.Range("L2").Value = 1.057961160651572
$ws.Range("M2").Value = 1.064261551352442
$ws.Range("N2").Value = 1.042483310182717

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03675138930629
$ws.Range("D3").Value = 1.041276423105893
$ws.Range("E3").Value = 1.056290218144739
$ws.Range("F3").Value = 1.062629952113701
$ws.Range("I3").Value = 1.041896282326701
$ws.Range("J3").Value = 1.041506877123633
$ws.Range("K3").Value = 1.04386912158745
$ws.Range("L3").Value = 1.058844018735362
$ws.Range("M3").Value = 1.065167683409268
$ws.Range("N3").Value = 1.042985936351657

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.037306567032181
$ws.Range("D4").Value = 1.041705044864208
$ws.Range("E4").Value = 1.056983750758079
$ws.Range("F4").Value = 1.063337593140697
$ws.Range("I4").Value = 1.042055275535106
$ws.Range("J4").Value = 1.041831174320949
$ws.Range("K4").Value = 1.04417466718927
$ws.Range("L4").Value = 1.059415782791319
$ws.Range("M4").Value = 1.065754334490204
$ws.Range("N4").Value = 1.04331069408819

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037540005795809
$ws.Range("D5").Value = 1.041885250708178
$ws.Range("E5").Value = 1.057275573749985
$ws.Range("F5").Value = 1.063635305090017
$ws.Range("I5").Value = 1.042121816256576
$ws.Range("J5").Value = 1.041967394415126
$ws.Range("K5").Value = 1.044302968750944
$ws.Range("L5").Value = 1.059656270506427
$ws.Range("M5").Value = 1.066001039255185
$ws.Range("N5").Value = 1.043447107630527

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037579203624648
$ws.Range("D6").Value = 1.041915508800035
$ws.Range("E6").Value = 1.057324587452258
$ws.Range("F6").Value = 1.063685305096704
$ws.Range("I6").Value = 1.042132971125861
$ws.Range("J6").Value = 1.041990259631148
$ws.Range("K6").Value = 1.0443245023352
$ws.Range("L6").Value = 1.059696656360478
$ws.Range("M6").Value = 1.066042466545715
$ws.Range("N6").Value = 1.043470005317777

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037309686088542
$ws.Range("D7").Value = 1.04170745273314
$ws.Range("E7").Value = 1.056987649082978
$ws.Range("F7").Value = 1.063341570320099
$ws.Range("I7").Value = 1.042056165835106
$ws.Range("J7").Value = 1.041832994951332
$ws.Range("K7").Value = 1.044176382150856
$ws.Range("L7").Value = 1.059418995736604
$ws.Range("M7").Value = 1.065757630671017
$ws.Range("N7").Value = 1.043312517304078

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036183506673539
$ws.Range("D8").Value = 1.040837928532484
$ws.Range("E8").Value = 1.055581526048281
$ws.Range("F8").Value = 1.061906681257882
$ws.Range("I8").Value = 1.04173258762162
$ws.Range("J8").Value = 1.041174685388179
$ws.Range("K8").Value = 1.043555995131528
$ws.Range("L8").Value = 1.058259423645738
$ws.Range("M8").Value = 1.064567715619206
$ws.Range("N8").Value = 1.042653272865836

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034200464816581
$ws.Range("D9").Value = 1.039306243540477
$ws.Range("E9").Value = 1.053112185626694
$ws.Range("F9").Value = 1.059385305089861
$ws.Range("I9").Value = 1.041152848176741
$ws.Range("J9").Value = 1.040011066301866
$ws.Range("K9").Value = 1.042458073514129
$ws.Range("L9").Value = 1.056219933447823
$ws.Range("M9").Value = 1.062473437769301
$ws.Range("N9").Value = 1.041488001306952

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032879463659069
$ws.Range("D10").Value = 1.038285555323419
$ws.Range("E10").Value = 1.051471662627645
$ws.Range("F10").Value = 1.057709189688787
$ws.Range("I10").Value = 1.040759966638777
$ws.Range("J10").Value = 1.039232956369811
$ws.Range("K10").Value = 1.0417230077041
$ws.Range("L10").Value = 1.054862878497533
$ws.Range("M10").Value = 1.061078975777589
$ws.Range("N10").Value = 1.040708786369543

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.032307712458052
$ws.Range("D11").Value = 1.037843705981529
$ws.Range("E11").Value = 1.050762658290792
$ws.Range("F11").Value = 1.056984560194035
$ws.Range("I11").Value = 1.040588336398212
$ws.Range("J11").Value = 1.038895476261042
$ws.Range("K11").Value = 1.041403988042027
$ws.Range("L11").Value = 1.05427588334715
$ws.Range("M11").Value = 1.060475574858138
$ws.Range("N11").Value = 1.04037082700031

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032095377343644
$ws.Range("D12").Value = 1.037679602121863
$ws.Range("E12").Value = 1.050499506303554
$ws.Range("F12").Value = 1.056715572312763
$ws.Range("I12").Value = 1.040524358965857
$ws.Range("J12").Value = 1.038770038814777
$ws.Range("K12").Value = 1.041285380955425
$ws.Range("L12").Value = 1.054057940629202
$ws.Range("M12").Value = 1.060251507099037
$ws.Range("N12").Value = 1.040245211418483

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032140922176963
$ws.Range("D13").Value = 1.037714802091243
$ws.Range("E13").Value = 1.050555944067575
$ws.Range("F13").Value = 1.056773263343824
$ws.Range("I13").Value = 1.040538092575132
$ws.Range("J13").Value = 1.038796949284066
$ws.Range("K13").Value = 1.041310827495369
$ws.Range("L13").Value = 1.054104685846962
$ws.Range("M13").Value = 1.060299567569822
$ws.Range("N13").Value = 1.040272160103725

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.032290159975053
$ws.Range("D14").Value = 1.037830140717013
$ws.Range("E14").Value = 1.050740901905074
$ws.Range("F14").Value = 1.056962322064307
$ws.Range("I14").Value = 1.040583052617454
$ws.Range("J14").Value = 1.03888510923685
$ws.Range("K14").Value = 1.041394186161245
$ws.Range("L14").Value = 1.054257866225628
$ws.Range("M14").Value = 1.060457052074961
$ws.Range("N14").Value = 1.040360445253754

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.032382115471936
$ws.Range("D15").Value = 1.037901207163061
$ws.Range("E15").Value = 1.050854887553407
$ws.Range("F15").Value = 1.057078830161178
$ws.Range("I15").Value = 1.040610724009863
$ws.Range("J15").Value = 1.038939416611699
$ws.Range("K15").Value = 1.041445531772276
$ws.Range("L15").Value = 1.054352258105596
$ws.Range("M15").Value = 1.060554091736403
$ws.Range("N15").Value = 1.040414829751304

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032917414279479
$ws.Range("D16").Value = 1.038314881953204
$ws.Range("E16").Value = 1.051518745460312
$ws.Range("F16").Value = 1.057757305019473
$ws.Range("I16").Value = 1.040771325405279
$ws.Range("J16").Value = 1.039255342225119
$ws.Range("K16").Value = 1.041744164656591
$ws.Range("L16").Value = 1.054901848506397
$ws.Range("M16").Value = 1.061119030202296
$ws.Range("N16").Value = 1.040731204015332

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.033253261007603
$ws.Range("D17").Value = 1.038574401048393
$ws.Range("E17").Value = 1.051935528579023
$ws.Range("F17").Value = 1.058183199862967
$ws.Range("I17").Value = 1.040871662538848
$ws.Range("J17").Value = 1.039453366448367
$ws.Range("K17").Value = 1.041931294011228
$ws.Range("L17").Value = 1.055246757998029
$ws.Range("M17").Value = 1.061473511427402
$ws.Range("N17").Value = 1.040929509455691

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033449178872685
$ws.Range("D18").Value = 1.038725785190963
$ws.Range("E18").Value = 1.052178761624925
$ws.Range("F18").Value = 1.058431727029294
$ws.Range("I18").Value = 1.040930041768007
$ws.Range("J18").Value = 1.039568817179286
$ws.Range("K18").Value = 1.042040372776668
$ws.Range("L18").Value = 1.055447997522408
$ws.Range("M18").Value = 1.06168031388109
$ws.Range("N18").Value = 1.04104512413989

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03351598586687
$ws.Range("D19").Value = 1.038777405104603
$ws.Range("E19").Value = 1.052261719927703
$ws.Range("F19").Value = 1.058516486995675
$ws.Range("I19").Value = 1.040949922843114
$ws.Range("J19").Value = 1.039608173804002
$ws.Range("K19").Value = 1.042077553811112
$ws.Range("L19").Value = 1.055516625103457
$ws.Range("M19").Value = 1.061750834865298
$ws.Range("N19").Value = 1.041084536655528

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.033217225326633
$ws.Range("D20").Value = 1.038546555953524
$ws.Range("E20").Value = 1.051890798192258
$ws.Range("F20").Value = 1.058137494004198
$ws.Range("I20").Value = 1.040860912382977
$ws.Range("J20").Value = 1.039432125845669
$ws.Range("K20").Value = 1.04191122410066
$ws.Range("L20").Value = 1.055209746300029
$ws.Range("M20").Value = 1.061435474851516
$ws.Range("N20").Value = 1.040908238688901

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032246212085536
$ws.Range("D21").Value = 1.037796175840766
$ws.Range("E21").Value = 1.050686430801969
$ws.Range("F21").Value = 1.056906644236751
$ws.Range("I21").Value = 1.040569819248824
$ws.Range("J21").Value = 1.038859150593506
$ws.Range("K21").Value = 1.041369642113018
$ws.Range("L21").Value = 1.054212755839391
$ws.Range("M21").Value = 1.060410675091587
$ws.Range("N21").Value = 1.040334449746159

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.031635922162055
$ws.Range("D22").Value = 1.037324490690954
$ws.Range("E22").Value = 1.049930376564757
$ws.Range("F22").Value = 1.05613375427469
$ws.Range("I22").Value = 1.040385488293526
$ws.Range("J22").Value = 1.0384984225161
$ws.Range("K22").Value = 1.041028497940031
$ws.Range("L22").Value = 1.053586449088903
$ws.Range("M22").Value = 1.059766703011515
$ws.Range("N22").Value = 1.039973209393506

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031959426761049
$ws.Range("D23").Value = 1.037574529177496
$ws.Range("E23").Value = 1.050331063255422
$ws.Range("F23").Value = 1.056543383440691
$ws.Range("I23").Value = 1.040483329534591
$ws.Range("J23").Value = 1.038689696117693
$ws.Range("K23").Value = 1.041209404347548
$ws.Range("L23").Value = 1.053918414731682
$ws.Range("M23").Value = 1.060108050463342
$ws.Range("N23").Value = 1.040164754625554

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.033233508225176
$ws.Range("D24").Value = 1.038559137919325
$ws.Range("E24").Value = 1.051911009525379
$ws.Range("F24").Value = 1.058158146174763
$ws.Range("I24").Value = 1.040865770365976
$ws.Range("J24").Value = 1.039441723724949
$ws.Range("K24").Value = 1.04192029304763
$ws.Range("L24").Value = 1.055226470108864
$ws.Range("M24").Value = 1.061452661820667
$ws.Range("N24").Value = 1.040917850198271

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034712952475888
$ws.Range("D25").Value = 1.03970215000097
$ws.Range("E25").Value = 1.053749566715933
$ws.Range("F25").Value = 1.060036297874178
$ws.Range("I25").Value = 1.041303853279963
$ws.Range("J25").Value = 1.040312310197594
$ws.Range("K25").Value = 1.042742466222921
$ws.Range("L25").Value = 1.056746733888052
$ws.Range("M25").Value = 1.063014557958875
$ws.Range("N25").Value = 1.041789673003564
